# Effort Commitment up to date
#
# This workbook tracks a team's effort-commitment schedule. The sheet is
# protected, so it must be unprotected before any cell contents can change.
# The edit:
#   - rolls the commitment "Start Date" forward from 2017-10-23 to 2019-10-14
#     (this recalculates G15 and every date/derived formula that depends on it)
#   - swaps which sprint columns (M & P) are flagged "Term" vs "Non-term" in
#     the week-type header row (row 20)
#   - correspondingly swaps the "Minimum"/"None" effort markers in the M and P
#     columns for the first four sprint rows (24-27), keeping the per-row
#     totals the same but shifting which week they land on
#   - leaves the cursor/selection on M30, where the edit was last made
#   - re-saves the workbook (no longer protected), matching the state Excel
#     leaves the file in after such an edit/review session

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it so the cells below can be edited.
$ws.Unprotect()

# --- Update the commitment Start Date (Year/Day) -----------------------
# D15 = Year, E15 = Month (unchanged), F15 = Day
$ws.Range("D15").Value2 = 2019
$ws.Range("F15").Value2 = 14

# --- Swap the Term / Non-term flags on row 20 for columns M and P ------
$ws.Range("M20").Value2 = "Term"
$ws.Range("P20").Value2 = "Non-term"

# --- Swap the Minimum / None effort markers for sprint rows 24-27 ------
# (columns M and P) while leaving every other column/value untouched.
$ws.Range("M24:M27").Value2 = "Minimum"
$ws.Range("P24:P27").Value2 = "None"

# --- Restore the active selection to M30, then save ---------------------
$ws.Range("M30").Select()

$wb.Save()
